$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.7339303333333334
$ws.Range("N2").Value = 2.201791
$ws.Range("O2").Value = 0.03574007706012852
$ws.Range("P2").Value = 0.03574007706012852
$ws.Range("Q2").Value = 0.04039430232944444
$ws.Range("R2").Value = 0.363548720965
$ws.Range("S2").Value = 0.03574007706012852
$ws.Range("T2").Value = 0.03574007706012852

# Row 3
$ws.Range("O3").Value = 0.3842514532634088
$ws.Range("P3").Value = 0.3842514532634088
$ws.Range("S3").Value = 0.3842514532634088
$ws.Range("T3").Value = 0.3842514532634088

# Row 4
$ws.Range("M4").Value = 4.974008666666667
$ws.Range("N4").Value = 14.922026
$ws.Range("O4").Value = 0.2422184299659874
$ws.Range("P4").Value = 0.2422184299659874
$ws.Range("Q4").Value = 0.2737611469988889
$ws.Range("R4").Value = 2.46385032299
$ws.Range("S4").Value = 0.2422184299659874
$ws.Range("T4").Value = 0.2422184299659874

# Row 5
$ws.Range("M5").Value = 2.087648
$ws.Range("N5").Value = 6.262943999999999
$ws.Range("O5").Value = 0.1016618294757629
$ws.Range("P5").Value = 0.1016618294757629
$ws.Range("Q5").Value = 0.1149006665066667
$ws.Range("R5").Value = 1.03410599856
$ws.Range("S5").Value = 0.1016618294757629
$ws.Range("T5").Value = 0.1016618294757629

# Row 6
$ws.Range("M6").Value = 4.848944666666667
$ws.Range("N6").Value = 14.546834
$ws.Range("O6").Value = 0.2361282102347124
$ws.Range("P6").Value = 0.2361282102347124
$ws.Range("Q6").Value = 0.2668778328788889
$ws.Range("R6").Value = 2.40190049591
$ws.Range("S6").Value = 0.2361282102347124
$ws.Range("T6").Value = 0.2361282102347124
